$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme
$cs.Colors(1).RGB  = 0
$cs.Colors(2).RGB  = 16777215
$cs.Colors(3).RGB  = 6968388
$cs.Colors(4).RGB  = 15132391
$cs.Colors(5).RGB  = 13998939
$cs.Colors(6).RGB  = 3243501
$cs.Colors(7).RGB  = 10855845
$cs.Colors(8).RGB  = 49407
$cs.Colors(9).RGB  = 12874308
$cs.Colors(10).RGB = 4697456
$cs.Colors(11).RGB = 12673797
$cs.Colors(12).RGB = 7491477
Write-Output "done"
